$p = $ppt.ActivePresentation

# --- Slide 20: "Zadanie: " -> "Zadanie 6: " ---
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(2)
$run20 = $sh20.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$run20.Text = "Zadanie 6: "

# --- Slide 23: fix typo "działata" -> "działała" and merge the three runs
#     of paragraph 5 into a single run (matching the author's re-typed text). ---
$s23 = $p.Slides.Item(23)
$sh23 = $s23.Shapes.Item(2)
$tr23 = $sh23.TextFrame.TextRange
$para5 = $tr23.Paragraphs(5, 1)
$sub = $tr23.Characters($para5.Start, $para5.Length)
$sub.Text = "Zablokowanie możliwości komunikacji do samej maszyny na której działała aplikacja i innych hostów w obrębie tej infrastruktury"
